$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Settings")
$ws2 = $wb.Worksheets.Item("LocalizedSettings")

# Rename Folder-related keys to OU-related keys (Settings -> LocalizedSettings sheet)
$ws2.Range("A27").Value = "GetOUsFailure"
$ws2.Range("A29").Value = "ChooseOUMessage"
$ws2.Range("A30").Value = "ChooseOUTitle"
$ws2.Range("A41").Value = "GetOUFailure"
$ws2.Range("A48").Value = "OUNotFound"

# Add new RobotNotFound row
$ws2.Range("A51").Value = "RobotNotFound"
$ws2.Range("B51").Value = "The robot named {0} was not found."
$ws2.Range("C51").Value = "{0}というロボットが見つかりませんでした。"

# Make LocalizedSettings the active/selected tab, with new selection
$ws2.Activate()
$ws2.Range("B29").Select()
